# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / handoff / handback timestamp
# cells with freshly generated values.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first file (row 2)
$wsOverview.Range("G2").Value = "2016-09-06 17:48:23"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime (row 2)
$wsZhCn.Range("H2").Value = "2016-09-06 17:48:16"
$wsZhCn.Range("K2").Value = "2016-09-06 17:48:46"

# de-de sheet: Correspond Handoff Datetime (row 2) and Correspond Handback DateTime (row 2)
$wsDeDe.Range("H2").Value = "2016-09-06 17:48:23"
$wsDeDe.Range("K2").Value = "2016-09-06 17:48:54"
